$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3.229934581839643, 0.4056676705253324, 0.01145294848073064, 0.05052024881177442, 3.628100717984239, 0.07973214163530429, 2.057488557387394, 0.1112702309759559, 0.4919442578564315, 2.150983187293811),
    @(3, 3.117137602727496, 0.3752551547852363, 0.01110465608546107, 0.05058848726883824, 3.615467914413358, 0.07973214163530429, 2.058799223821048, 0.1119254443326003, 0.4853832162078078, 2.17463226720351),
    @(4, 3.049638617032031, 0.3568102768786048, 0.01088656952543943, 0.0506343506248248, 3.609672594614665, 0.07973214163530429, 2.060670984191844, 0.1123519967546915, 0.4815950824088304, 2.189888701198651),
    @(5, 3.022573837985931, 0.3493507527682596, 0.01079659825546386, 0.05065403867706086, 3.607802938907923, 0.07973214163530429, 2.061701419700398, 0.1125319266524247, 0.4801118100636188, 2.19629057731856),
    @(6, 3.018106405460458, 0.3481155259416084, 0.01078159119555977, 0.05065736820472655, 3.60752216412223, 0.07973214163530429, 2.061888671596051, 0.1125621729404136, 0.4798691626971703, 2.197364754984271),
    @(7, 3.049271824786217, 0.3567094454075743, 0.01088536063736179, 0.05063461210027653, 3.609645389360224, 0.07973214163530429, 2.060683798037161, 0.1123543986162749, 0.4815748338533581, 2.189974291445662),
    @(8, 3.19067668873538, 0.3951336593515578, 0.01133371497406799, 0.05054295539598419, 3.623337128173162, 0.07973214163530429, 2.05771863206455, 0.1114911210119161, 0.4896321047580727, 2.158984514163468),
    @(9, 3.481976509065532, 0.472326026168048, 0.01218093158036382, 0.05039461646826149, 3.665808928162733, 0.07973214163530429, 2.06040110864538, 0.1099902342023151, 0.5073421517624865, 2.104063065235664),
    @(10, 3.704630519717625, 0.5302095191993317, 0.0127861024474214, 0.05030470144833332, 3.706630235277146, 0.07973214163530429, 2.067598579541382, 0.1090040058790791, 0.5215237486767563, 2.067292283782372),
    @(11, 3.807820368240073, 0.5568069704131631, 0.01305814622214818, 0.05026792373429656, 3.727310351989217, 0.07973214163530429, 2.072018281737897, 0.1085805149330152, 0.5282308357394214, 2.051344766703167),
    @(12, 3.847170789188965, 0.5669177241858847, 0.01316073416475483, 0.050254589131697, 3.735446427469299, 0.07973214163530429, 2.07385745002675, 0.1084237577124636, 0.5308075008476294, 2.045418295348888),
    @(13, 3.838683740471424, 0.5647384532508113, 0.01313865848409712, 0.05025743464650101, 3.733680590397796, 0.07973214163530429, 2.073453977717776, 0.1084573577532311, 0.5302509307158516, 2.046689656927995),
    @(14, 3.811052240728827, 0.5576380054159245, 0.01306659457172188, 0.05026681482218764, 3.727973591158616, 0.07973214163530429, 2.072166269284267, 0.1085675461188771, 0.5284420809515638, 2.050854935273961),
    @(15, 3.794162947888083, 0.5532938578675726, 0.01302239860828891, 0.05027263756374278, 3.724517649825543, 0.07973214163530429, 2.071399090528473, 0.1086355095515064, 0.5273389070943324, 2.053420951952695),
    @(16, 3.697925202751264, 0.5284767187322927, 0.01276826159379851, 0.05030718789041122, 3.705321333699061, 0.07973214163530429, 2.067332860482892, 0.1090321874899294, 0.521090576614867, 2.068350213968827),
    @(17, 3.639374640395829, 0.513320772219231, 0.01261155118413093, 0.05032943929370781, 3.694086576233929, 0.07973214163530429, 2.065132321629775, 0.1092819729943315, 0.5173229823735852, 2.077708804767298),
    @(18, 3.60587700475827, 0.5046284676064374, 0.01252110713173415, 0.05034262603158646, 3.68782318057967, 0.07973214163530429, 2.063974404798159, 0.109428010579399, 0.5151800362677221, 2.083164983633665),
    @(19, 3.594566012996438, 0.5016896764702778, 0.01249043020332863, 0.05034715755602526, 3.685736556365526, 0.07973214163530429, 2.06360084044546, 0.1094778632382107, 0.5144586048297271, 2.085024940853714),
    @(20, 3.645588904853298, 0.5149315571762827, 0.01262826491254998, 0.05032703041080988, 3.695261976429805, 0.07973214163530429, 2.065355412732416, 0.1092551379165307, 0.5177215572349496, 2.076704971025272),
    @(21, 3.819160819835645, 0.559722518883234, 0.01308777282858031, 0.05026404356832992, 3.729641586396212, 0.07973214163530429, 2.072540001804128, 0.1085350832147149, 0.5289723837990437, 2.049628435505909),
    @(22, 3.934201533063003, 0.589222903301561, 0.01338560466940386, 0.05022633016888367, 3.753888931911035, 0.07973214163530429, 2.078200720091829, 0.1080855223222486, 0.5365401895164297, 2.032588318796094),
    @(23, 3.872655280549907, 0.5734570198055735, 0.0132268606033854, 0.05024614291512131, 3.740784449413354, 0.07973214163530429, 2.07509090801986, 0.1083235388163839, 0.5324814421737045, 2.041622798559317),
    @(24, 3.642778925567882, 0.5142032557756693, 0.01262070972596341, 0.05032811823916128, 3.694729968814443, 0.07973214163530429, 2.065254219393836, 0.1092672624773332, 0.5175412896109037, 2.077158568188288),
    @(25, 3.40166285222341, 0.451241366429997, 0.0119549342124099, 0.05043139226777416, 3.652637135319651, 0.07973214163530429, 2.058760982826868, 0.1103757688884182, 0.5023460782418141, 2.118293131375566)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]   # B
    $ws.Cells.Item($row, 3).Value = $entry[2]   # C
    $ws.Cells.Item($row, 4).Value = $entry[3]   # D
    $ws.Cells.Item($row, 5).Value = $entry[4]   # E
    $ws.Cells.Item($row, 6).Value = $entry[5]   # F
    $ws.Cells.Item($row, 8).Value = $entry[6]   # H
    $ws.Cells.Item($row, 9).Value = $entry[7]   # I
    $ws.Cells.Item($row, 10).Value = $entry[8]  # J
    $ws.Cells.Item($row, 12).Value = $entry[9]  # L
    $ws.Cells.Item($row, 14).Value = $entry[10] # N
}
